# Generate Report for Handback
# Rewrites the two request-id/hash pairs that the handback CI run produced
# this time around, across the Overview, zh-cn and de-de sheets, and updates
# the display text of every hyperlink that showed the old names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values for this handback run
# ---------------------------------------------------------------------
$oldMd1 = "1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md"
$newMd1 = "cc9c4d88-a1e4-4904-acb4-021401f05c23.md"

$oldMd2 = "896cbfaa-c4b3-4761-9744-ecdb62318b5b.md"
$newMd2 = "ffff76ca6570-403e-4963-9d34-7be281117591.md"

$newXlfZh = "cc9c4d88-a1e4-4904-acb4-021401f05c23.6c0a929d12438973a48157a11a42e7268d01a887.zh-cn.xlf"
$newXlfDe = "cc9c4d88-a1e4-4904-acb4-021401f05c23.6c0a929d12438973a48157a11a42e7268d01a887.de-de.xlf"

$newHandoffZh1 = "2016-03-17 14:50:52"
$newHandbackZh1 = "2016-03-17 14:51:10"
$newHandoffDe1 = "2016-03-17 14:50:56"
$newHandbackDe1 = "2016-03-17 14:51:16"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = $newMd1
$ws1.Range("A3").Value = $newMd2

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", $newMd1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", $newMd2)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newMd1
$ws2.Range("D2").Value = $newXlfZh
$ws2.Range("E2").Value = $newHandoffZh1
$ws2.Range("F2").Value = $newMd1
$ws2.Range("G2").Value = $newXlfZh
$ws2.Range("H2").Value = $newHandbackZh1

$ws2.Range("A3").Value = $newMd2
$ws2.Range("D3").Value = $newXlfZh
$ws2.Range("E3").Value = $newHandoffZh1
$ws2.Range("F3").Value = $newMd2
$ws2.Range("G3").Value = $newXlfZh
$ws2.Range("H3").Value = $newHandbackZh1

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", $newMd1)
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8abbb53f2c64c021fa370a9a834f5a55f3b6eba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.bd939de73871ba63faeab632f3b9f3b5a97e35f4.zh-cn.xlf", "", "", $newXlfZh)
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fbd58c8bf64c8bcf95dac94e5597a4fa9c349ba8/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", $newMd1)
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/67491639c036df42a448bbdabf11813d2a0c0288/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.bd939de73871ba63faeab632f3b9f3b5a97e35f4.zh-cn.xlf", "", "", $newXlfZh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", $newMd2)
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8abbb53f2c64c021fa370a9a834f5a55f3b6eba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/896cbfaa-c4b3-4761-9744-ecdb62318b5b.f7bba4fba79bbad95d1f757e5d35650020e9fff3.zh-cn.xlf", "", "", $newXlfZh)
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fbd58c8bf64c8bcf95dac94e5597a4fa9c349ba8/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", $newMd2)
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/67491639c036df42a448bbdabf11813d2a0c0288/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/896cbfaa-c4b3-4761-9744-ecdb62318b5b.f7bba4fba79bbad95d1f757e5d35650020e9fff3.zh-cn.xlf", "", "", $newXlfZh)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newMd1
$ws3.Range("D2").Value = $newXlfDe
$ws3.Range("E2").Value = $newHandoffDe1
$ws3.Range("F2").Value = $newMd1
$ws3.Range("G2").Value = $newXlfDe
$ws3.Range("H2").Value = $newHandbackDe1

$ws3.Range("A3").Value = $newMd2
$ws3.Range("D3").Value = $newXlfDe
$ws3.Range("E3").Value = $newHandoffDe1
$ws3.Range("F3").Value = $newMd2
$ws3.Range("G3").Value = $newXlfDe
$ws3.Range("H3").Value = $newHandbackDe1

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", $newMd1)
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bcd14aa45121d51d5299f46177d81b756362cbe2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.bd939de73871ba63faeab632f3b9f3b5a97e35f4.de-de.xlf", "", "", $newXlfDe)
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/9793281fe2f2606253fcf6f2a8fade70b4e7f002/e2e/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.md", "", "", $newMd1)
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ea585580f355bdbd2ef14daad10bc190a879a887/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1c0d43dd-f7e9-4bf9-abd4-f16f449c04bf.bd939de73871ba63faeab632f3b9f3b5a97e35f4.de-de.xlf", "", "", $newXlfDe)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", $newMd2)
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d4731bbef80a452a3fa136653c559bc7ba8d16a5/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bcd14aa45121d51d5299f46177d81b756362cbe2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/896cbfaa-c4b3-4761-9744-ecdb62318b5b.f7bba4fba79bbad95d1f757e5d35650020e9fff3.de-de.xlf", "", "", $newXlfDe)
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/9793281fe2f2606253fcf6f2a8fade70b4e7f002/e2e/896cbfaa-c4b3-4761-9744-ecdb62318b5b.md", "", "", $newMd2)
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ea585580f355bdbd2ef14daad10bc190a879a887/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/896cbfaa-c4b3-4761-9744-ecdb62318b5b.f7bba4fba79bbad95d1f757e5d35650020e9fff3.de-de.xlf", "", "", $newXlfDe)
